# Notes test-data sheet: keep only the "Note" column (previously column C),
# dropping "List name" (A), "Task title" (B) and the unused formatting-only
# column (D). Deleting D, then B, then A (right-to-left) shifts the "Note"
# column all the way down to column A while preserving its original style
# and column width.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D:D").Delete()
$ws.Range("B:B").Delete()
$ws.Range("A:A").Delete()
